$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.985.17'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '1.741.23'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.95'
$ws.Range('E5').Value = '  +2.85%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5044'
$ws.Range('E7').Value = '  -4.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2738'
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06174'
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = '1.752.42'
$ws.Range('E10').Value = '  +0.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07252'
$ws.Range('E11').Value = '  +0.84%  '
$ws.Range('E12').Value = '  +1.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.09'
$ws.Range('E13').Value = '  -0.37%  '
$ws.Range('E14').Value = '  +1.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.55'
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('D18').Value = '26.005.48'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.90'
$ws.Range('E19').Value = '  +0.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006858'
$ws.Range('E20').Value = '  +1.10%  '
$ws.Range('D21').Value = '1.977.16'
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.473'
$ws.Range('E22').Value = '  +2.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.713'
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.368'
$ws.Range('E24').Value = '  +2.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '135.37'
$ws.Range('E25').Value = '  -3.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.502'
$ws.Range('E26').Value = '  -0.85%  '
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.782'
$ws.Range('E28').Value = '  +0.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '105.37'
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.913'
$ws.Range('E30').Value = '  +1.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08159'
$ws.Range('E31').Value = '  -3.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.672'
$ws.Range('E32').Value = '  +0.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04682'
$ws.Range('E33').Value = '  +1.89%  '
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9961'
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6128'
$ws.Range('E36').Value = '  -1.86%  '
$ws.Range('E37').Value = '  +2.34%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01620'
$ws.Range('E38').Value = '  +0.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.926'
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.000'
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '100.72'
$ws.Range('E41').Value = '  +1.85%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.3913'
$ws.Range('E42').Value = '  +0.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7626'
$ws.Range('E43').Value = '  +1.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.007'
$ws.Range('E44').Value = '  +1.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1158'
$ws.Range('E45').Value = '  +1.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.320'
$ws.Range('E46').Value = '  +1.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.52'
$ws.Range('E47').Value = '  +1.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05300'
$ws.Range('E49').Value = '  -0.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3468'
$ws.Range('E50').Value = '  +0.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.588'
